$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'321.72"
$ws.Range("D3").Value = "'42.87"
$ws.Range("E3").Value = "'-6.58%"
$ws.Range("D4").Value = "'5.174"
$ws.Range("E4").Value = "'-8.30%"
$ws.Range("D5").Value = "'0.08154"
$ws.Range("E5").Value = "'-2.65%"
$ws.Range("D6").Value = "'4.328"
$ws.Range("E6").Value = "'-3.51%"
$ws.Range("D7").Value = "'1.813"
$ws.Range("E7").Value = "'-12.29%"
$ws.Range("D8").Value = "'0.9452"
$ws.Range("E8").Value = "'-4.33%"
$ws.Range("D9").Value = "'0.1118"
$ws.Range("E9").Value = "'-3.26%"
$ws.Range("D10").Value = "'0.1850"
$ws.Range("E10").Value = "'-4.14%"
$ws.Range("D11").Value = "'0.09339"
$ws.Range("E11").Value = "'-6.30%"
$ws.Range("D12").Value = "'0.04622"
$ws.Range("E12").Value = "'-1.13%"
$ws.Range("D13").Value = "'7.411"
$ws.Range("E13").Value = "'-28.75%"
$ws.Range("E14").Value = "'0.04%"
$ws.Range("D15").Value = "'0.001294"
$ws.Range("E15").Value = "'1.40%"
$ws.Range("D16").Value = "'0.005642"
$ws.Range("E16").Value = "'-7.93%"
$ws.Range("D17").Value = "'3.360"
$ws.Range("E17").Value = "'-0.50%"
$ws.Range("D18").Value = "'2.504"
$ws.Range("E18").Value = "'-3.05%"
$ws.Range("D19").Value = "'0.3366"
$ws.Range("E19").Value = "'0.00%"
$ws.Range("D20").Value = "'0.1389"
$ws.Range("E20").Value = "'-0.96%"
$ws.Range("D21").Value = "'0.2622"
$ws.Range("E21").Value = "'-1.29%"
$ws.Range("D22").Value = "'0.04180"
$ws.Range("E22").Value = "'-0.72%"
$ws.Range("D23").Value = "'0.001250"
$ws.Range("E23").Value = "'-5.02%"
$ws.Range("D24").Value = "'0.004295"
$ws.Range("E24").Value = "'-7.61%"
$ws.Range("E25").Value = "'-13.50%"
$ws.Range("D26").Value = "'0.0002982"
$ws.Range("E26").Value = "'-20.50%"
$ws.Range("D38").Value = "'0.02699"
$ws.Range("E38").Value = "'-3.08%"
$ws.Range("D39").Value = "'0.05526"
$ws.Range("E39").Value = "'-3.73%"
$ws.Range("D40").Value = "'0.007965"
$ws.Range("E40").Value = "'2.53%"
$ws.Range("D41").Value = "'0.1394"
$ws.Range("E41").Value = "'-2.94%"
$ws.Range("D42").Value = "'0.006555"
$ws.Range("E42").Value = "'-10.07%"
$ws.Range("D43").Value = "'0.002121"
$ws.Range("E43").Value = "'0.03%"
$ws.Range("D44").Value = "'0.007466"
$ws.Range("E44").Value = "'-17.67%"
$ws.Range("D45").Value = "'0.3203"
$ws.Range("E45").Value = "'-6.06%"
$ws.Range("D46").Value = "'0.00006988"
$ws.Range("E46").Value = "'-1.74%"
$ws.Range("E47").Value = "'-0.19%"
$ws.Range("D48").Value = "'0.003468"
$ws.Range("E48").Value = "'-0.87%"
$ws.Range("D49").Value = "'0.003533"
$ws.Range("E49").Value = "'0.72%"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("E50").Value = "'-0.19%"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("E51").Value = "'-0.19%"

# Reset style on touched cells so the quote-prefix flag introduced by the
# leading apostrophe (used to force text entry) does not linger as a style delta.
$ws.Range("D2:E51").Style = "Normal"

